$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 66
$ws1.Cells.Item(3, 6).Value = 312
$ws1.Cells.Item(4, 6).Value = 167
$ws1.Cells.Item(5, 6).Value = 196
$ws1.Cells.Item(6, 6).Value = 333
$ws1.Cells.Item(7, 7).Value = 69
$ws1.Cells.Item(8, 6).Value = 2169
$ws1.Cells.Item(9, 6).Value = 373
$ws1.Cells.Item(10, 6).Value = 5329

# Sheet "全部类型" (sheet index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 66
$ws4.Cells.Item(4, 6).Value = 312
$ws4.Cells.Item(5, 6).Value = 167
$ws4.Cells.Item(6, 6).Value = 196
$ws4.Cells.Item(7, 6).Value = 333
$ws4.Cells.Item(8, 7).Value = 69
$ws4.Cells.Item(11, 6).Value = 2169
$ws4.Cells.Item(12, 6).Value = 373
$ws4.Cells.Item(13, 6).Value = 5329

$wb.Save()
